$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "表头1"
$ws.Range("B1").Value = "表头2"
$ws.Range("C1").Value = "表头3"

# Row 2
$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = "B"
$ws.Range("C2").Value = "C"

# Row 3
$ws.Range("A3").Value = "D"
$ws.Range("B3").Value = "E"
$ws.Range("C3").Value = "F"

$ws.Range("C1").Select()
